$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format so numeric-looking strings are preserved as text (matches original inlineStr cells)
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "290.07"
$ws.Range("E2").Value = "-3.81%"
$ws.Range("G2").Value = "23"

$ws.Range("E3").Value = "-3.82%"
$ws.Range("G3").Value = "23"

$ws.Range("E4").Value = "-1.99%"
$ws.Range("G4").Value = "23"

$ws.Range("D5").Value = "0.07162"
$ws.Range("E5").Value = "-9.23%"
$ws.Range("G5").Value = "23"

$ws.Range("D6").Value = "1.839"
$ws.Range("E6").Value = "-13.69%"
$ws.Range("G6").Value = "23"

$ws.Range("D7").Value = "7.647"
$ws.Range("E7").Value = "-2.03%"
$ws.Range("G7").Value = "23"

$ws.Range("D8").Value = "3.776"
$ws.Range("E8").Value = "-1.71%"
$ws.Range("G8").Value = "23"

$ws.Range("D9").Value = "0.8943"
$ws.Range("E9").Value = "-3.42%"
$ws.Range("G9").Value = "23"

$ws.Range("D10").Value = "0.1646"
$ws.Range("E10").Value = "-5.59%"
$ws.Range("G10").Value = "23"

$ws.Range("D11").Value = "0.07574"
$ws.Range("E11").Value = "-4.67%"
$ws.Range("G11").Value = "23"

$ws.Range("D12").Value = "0.08049"
$ws.Range("E12").Value = "-6.81%"
$ws.Range("G12").Value = "23"

$ws.Range("D13").Value = "0.02982"
$ws.Range("E13").Value = "-3.67%"
$ws.Range("G13").Value = "23"

$ws.Range("D14").Value = "0.1001"
$ws.Range("E14").Value = "-0.07%"
$ws.Range("G14").Value = "23"

$ws.Range("D15").Value = "0.001493"
$ws.Range("E15").Value = "-1.07%"
$ws.Range("G15").Value = "23"

$ws.Range("D16").Value = "0.005746"
$ws.Range("E16").Value = "-3.14%"
$ws.Range("G16").Value = "23"

$ws.Range("G17").Value = "23"

$ws.Range("D18").Value = "3.471"
$ws.Range("E18").Value = "0.28%"
$ws.Range("G18").Value = "23"

$ws.Range("D19").Value = "2.108"
$ws.Range("E19").Value = "-6.56%"
$ws.Range("G19").Value = "23"

$ws.Range("D20").Value = "0.3277"
$ws.Range("E20").Value = "-0.31%"
$ws.Range("G20").Value = "23"

$ws.Range("E21").Value = "-1.36%"
$ws.Range("G21").Value = "23"

$ws.Range("D22").Value = "4.264"
$ws.Range("E22").Value = "0.03%"
$ws.Range("G22").Value = "23"

$ws.Range("D23").Value = "0.2001"
$ws.Range("E23").Value = "11.56%"
$ws.Range("G23").Value = "23"

$ws.Range("D24").Value = "0.04475"
$ws.Range("E24").Value = "-2.90%"
$ws.Range("G24").Value = "23"

$ws.Range("D25").Value = "0.001213"
$ws.Range("E25").Value = "-1.97%"
$ws.Range("G25").Value = "23"

$ws.Range("D26").Value = "0.004668"
$ws.Range("E26").Value = "4.96%"
$ws.Range("G26").Value = "23"

$ws.Range("D27").Value = "0.0001251"
$ws.Range("E27").Value = "0.10%"
$ws.Range("G27").Value = "23"

$ws.Range("G28").Value = "23"

$ws.Range("G29").Value = "23"

$ws.Range("G30").Value = "23"

$ws.Range("G31").Value = "23"

$ws.Range("G32").Value = "23"

$ws.Range("G33").Value = "23"

$ws.Range("G34").Value = "23"

$ws.Range("G35").Value = "23"

$ws.Range("G36").Value = "23"

$ws.Range("G37").Value = "23"

$ws.Range("G38").Value = "23"

$ws.Range("D39").Value = "0.01642"
$ws.Range("E39").Value = "-4.32%"
$ws.Range("G39").Value = "23"

$ws.Range("D40").Value = "0.04354"
$ws.Range("E40").Value = "-8.91%"
$ws.Range("G40").Value = "23"

$ws.Range("D41").Value = "0.007380"
$ws.Range("E41").Value = "-1.06%"
$ws.Range("G41").Value = "23"

$ws.Range("D42").Value = "0.1307"
$ws.Range("E42").Value = "-3.65%"
$ws.Range("G42").Value = "23"

$ws.Range("E43").Value = "-14.72%"
$ws.Range("G43").Value = "23"

$ws.Range("D44").Value = "0.01026"
$ws.Range("E44").Value = "-8.75%"
$ws.Range("G44").Value = "23"

$ws.Range("D45").Value = "0.00005842"
$ws.Range("E45").Value = "-2.34%"
$ws.Range("G45").Value = "23"

$ws.Range("E46").Value = "0.10%"
$ws.Range("G46").Value = "23"

$ws.Range("D47").Value = "2.213"
$ws.Range("E47").Value = "168.74%"
$ws.Range("G47").Value = "23"

$ws.Range("E48").Value = "-11.36%"
$ws.Range("G48").Value = "23"

$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").Value = "0.10%"
$ws.Range("G49").Value = "23"

$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").Value = "0.10%"
$ws.Range("G50").Value = "23"

$ws.Range("G51").Value = "23"
